# Update mapa_interactivo_PEBCOM.xlsx - automatic map update
# 1) Insert a new record (row 38) for "José A. Cabrera 3086", shifting the
#    existing rows 38-71 down to 39-72.
# 2) Append a new record (row 73) for "ACEVEDO 310".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row at position 38 (pushes old rows 38.. down by one) ---
$ws.Rows.Item(38).Insert()

# Make sure text-like columns keep their text representation instead of
# being auto-converted to dates/numbers by Excel.
$ws.Range("A38:H38").NumberFormat = "@"
$ws.Range("J38:L38").NumberFormat = "@"
$ws.Range("O38:P38").NumberFormat = "@"

$ws.Range("A38").Value = "803825124"
$ws.Range("B38").Value = "3/7/2025"
$ws.Range("C38").Value = "José A. Cabrera 3086"
$ws.Range("D38").Value = "2"
$ws.Range("E38").Value = "803825124"
$ws.Range("F38").Value = "PEBCOM"
$ws.Range("G38").Value = "Pendiente"
$ws.Range("H38").Value = "Desmontar columna y transferir a comunitaria"
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = "Cambio"
$ws.Range("K38").Value = "Sin equipos"
$ws.Range("L38").Value = "Pasante"
$ws.Range("M38").Value = -58.41002
$ws.Range("N38").Value = -34.596998
$ws.Range("O38").Value = "Almagro"
$ws.Range("P38").Value = "Capital Sur"

# --- Append a brand-new row 73 at the bottom of the table ---
$ws.Range("A73:H73").NumberFormat = "@"
$ws.Range("J73:L73").NumberFormat = "@"
$ws.Range("O73:P73").NumberFormat = "@"

$ws.Range("A73").Value = "6471"
$ws.Range("B73").Value = "7/25/2025"
$ws.Range("C73").Value = "ACEVEDO 310"
$ws.Range("D73").Value = "15"
$ws.Range("E73").Value = "808533124"
$ws.Range("F73").Value = "PEBCOM"
$ws.Range("G73").Value = "Pendiente"
$ws.Range("H73").Value = "Picada"
$ws.Range("I73").Value = 1
$ws.Range("J73").Value = "Cambio"
$ws.Range("K73").Value = "Sin equipos"
$ws.Range("L73").Value = "Pasante"
$ws.Range("M73").Value = -58.44163
$ws.Range("N73").Value = -34.598788
$ws.Range("O73").Value = "Paternal"
$ws.Range("P73").Value = "Capital Norte"
